$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Septiembre de 2020 a las 01:08"

# Swap country names: Argentina now ranks above Chile (row 13/14)
$ws.Range("A13").Value = "Argentina"
$ws.Range("A14").Value = "Chile"

# Swap country names: Japon now ranks above Polonia (row 47/48)
$ws.Range("A47").Value = "Japon"
$ws.Range("A48").Value = "Polonia"

# Update numeric data cells
# Row 4
$ws.Range("B4").Value = 6208875
$ws.Range("C4").Value = 35639
$ws.Range("D4").Value = 3450243
$ws.Range("E4").Value = 2570955
$ws.Range("G4").Value = 453
$ws.Range("H4").Value = 187677

# Row 10
$ws.Range("B10").Value = 615168
$ws.Range("C10").Value = 7230
$ws.Range("D10").Value = 459475
$ws.Range("E10").Value = 136030
$ws.Range("G10").Value = 299
$ws.Range("H10").Value = 19663

# Row 13
$ws.Range("B13").Value = 417735
$ws.Range("C13").Value = 9309
$ws.Range("D13").Value = 301195
$ws.Range("E13").Value = 107880
$ws.Range("G13").Value = 203
$ws.Range("H13").Value = 8660

# Row 14
$ws.Range("B14").Value = 411726
$ws.Range("C14").Value = 1752
$ws.Range("D14").Value = 383879
$ws.Range("E14").Value = 16558
$ws.Range("G14").Value = 45
$ws.Range("H14").Value = 11289

# Row 23
$ws.Range("B23").Value = 244792
$ws.Range("C23").Value = 1497
$ws.Range("E23").Value = 17579

# Row 27
$ws.Range("B27").Value = 128948
$ws.Range("C27").Value = 1008
$ws.Range("D27").Value = 114227
$ws.Range("E27").Value = 5595
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = 9126

# Row 34
$ws.Range("B34").Value = 98939
$ws.Range("C34").Value = 212
$ws.Range("D34").Value = 72929
$ws.Range("E34").Value = 20589
$ws.Range("G34").Value = 22
$ws.Range("H34").Value = 5421

# Row 43
$ws.Range("B43").Value = 74074
$ws.Range("C43").Value = 162
$ws.Range("D43").Value = 62068
$ws.Range("E43").Value = 9246
$ws.Range("G43").Value = 20
$ws.Range("H43").Value = 2760

# Row 47
$ws.Range("B47").Value = 67865
$ws.Range("C47").Value = 601
$ws.Range("D47").Value = 56802
$ws.Range("E47").Value = 9784
$ws.Range("G47").Value = 15
$ws.Range("H47").Value = 1279

# Row 48
$ws.Range("B48").Value = 67372
$ws.Range("C48").Value = 502
$ws.Range("D48").Value = 46638
$ws.Range("E48").Value = 18695
$ws.Range("G48").Value = 6
$ws.Range("H48").Value = 2039

# Row 53
$ws.Range("B53").Value = 54008
$ws.Range("C53").Value = 143
$ws.Range("D53").Value = 41638
$ws.Range("E53").Value = 11357

# Row 55
$ws.Range("B55").Value = 51972
$ws.Range("C55").Value = 398
$ws.Range("D55").Value = 49047
$ws.Range("E55").Value = 2735

# Row 90
$ws.Range("B90").Value = 10782
$ws.Range("C90").Value = 139
$ws.Range("E90").Value = 1170

# Row 99
$ws.Range("B99").Value = 8533
$ws.Range("C99").Value = 28
$ws.Range("D99").Value = 7264
$ws.Range("E99").Value = 1216

# Row 106
$ws.Range("B106").Value = 6497
$ws.Range("C106").Value = 85
$ws.Range("D106").Value = 5221
$ws.Range("E106").Value = 1074
$ws.Range("G106").Value = 6
$ws.Range("H106").Value = 202

# Row 151
$ws.Range("B151").Value = 1595
$ws.Range("C151").Value = 10
$ws.Range("D151").Value = 1409
$ws.Range("E151").Value = 142

# Row 165
$ws.Range("B165").Value = 1007
$ws.Range("C165").Value = 14
$ws.Range("D165").Value = 173
$ws.Range("E165").Value = 821
